$wb = $excel.ActiveWorkbook

# The metadata table lives on the "Metadata" worksheet.
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$ws.Range("B3").Value = "0.0.0"

# Title: "Extension of Race" -> "Race"
$ws.Range("B5").Value = "Race"

# Date: 2023-11-21T19:08:35-03:00 -> 2024-01-16T20:08:55-03:00
$ws.Range("B8").Value = "2024-01-16T20:08:55-03:00"

# Description: updated wording
$ws.Range("B12").Value = "Extension that represents the race of a patient."

# The "Elements" sheet repeats the Title/Description of the root
# "Extension" element (row 2, columns L="Short" / M="Definition") using the
# very same shared-string text as Metadata!B5 / Metadata!B12. Update them
# too so both usages stay in sync, matching the shared string table edit.
$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("L2").Value = "Race"
$ws2.Range("M2").Value = "Extension that represents the race of a patient."
